$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Overview")

# --- Update the database: drop the oldest reported period (column D, "6 ماهه
# منتهی به 1399/06") and shift every later period one column to the left.
$ws.Columns("D:D").Delete()

# Copy the formatting of the (now last populated) column L into the freshly
# opened column M so the new period's cells keep the same number formats /
# styles as the rest of the table.
$ws.Range("L8:L27").Copy()
$ws.Range("M8:M27").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Add the new latest period, "12 ماهه منتهی به 1401/12", in column M.
$ws.Range("M8").Value = "12 ماهه منتهی به 1401/12"
$ws.Range("M9").Value = "'1402-02-28"

# The publish-date note for the prior column (now I9, "3 ماهه منتهی به 1401/03")
# was refreshed to reflect the latest disclosure.
$ws.Range("I9").Value = "1402-02-28 (8)"

# New column's financial figures (rial, cumulative).
$ws.Range("M11").Value = 13084540
$ws.Range("M12").Value = -6968838
$ws.Range("M13").Value = 6115702
$ws.Range("M14").Value = -596407
$ws.Range("M15").Value = 0
$ws.Range("M16").Value = 328860
$ws.Range("M17").Value = 5848155
$ws.Range("M18").Value = -461471
$ws.Range("M19").Value = 48199
$ws.Range("M20").Value = 5434883
$ws.Range("M21").Value = -757110
$ws.Range("M22").Value = 4677773
$ws.Range("M23").Value = 0
$ws.Range("M24").Value = 4677773
$ws.Range("M25").Value = 6542
$ws.Range("M26").Value = 715000
$ws.Range("M27").Value = 6542
